$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A145").Value = "PVB"
